$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Kode Mata Kuliah
$ws.Range("C8").Value = "EE0103-19"

# Nama Mata Kuliah
$ws.Range("C9").Value = "Matematika Diskret dan Logika"

# Pengampu
$ws.Range("C11").Value = "Dr.Ir. Augustinus Sujono M.T."

# Hidden/hash value stored merged in B12:C12
$ws.Range("B12").Value = "caFNPWvAHl//MJRM6J0jcw=="

# Tahun Ajaran numeric value
$ws.Range("C6").Value = 2016
